$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("M2").Value = 1.716657
$ws.Range("N2").Value = 5.149971
$ws.Range("O2").Value = 0.3840886036988016
$ws.Range("P2").Value = 0.3840886036988015
$ws.Range("Q2").Value = 0.6119292819429999
$ws.Range("R2").Value = 5.507363537487
$ws.Range("S2").Value = 0.3840886036988016
$ws.Range("T2").Value = 0.3840886036988015

# Row 3
$ws.Range("O3").Value = 0.07870146593648156
$ws.Range("P3").Value = 0.07870146593648154
$ws.Range("S3").Value = 0.07870146593648156
$ws.Range("T3").Value = 0.07870146593648154

# Row 4
$ws.Range("M4").Value = 1.677572333333333
$ws.Range("N4").Value = 5.032717
$ws.Range("O4").Value = 0.3753437146230962
$ws.Range("P4").Value = 0.3753437146230962
$ws.Range("Q4").Value = 0.5979969401832221
$ws.Range("R4").Value = 5.381972461648999
$ws.Range("S4").Value = 0.3753437146230962
$ws.Range("T4").Value = 0.3753437146230962

# Row 5
$ws.Range("M5").Value = 0.7234496666666667
$ws.Range("N5").Value = 2.170349
$ws.Range("O5").Value = 0.1618662157416207
$ws.Range("P5").Value = 0.1618662157416207
$ws.Range("Q5").Value = 0.2578849677281111
$ws.Range("R5").Value = 2.320964709552999
$ws.Range("S5").Value = 0.1618662157416207
$ws.Range("T5").Value = 0.1618662157416207
